$wb = $excel.ActiveWorkbook

# --- Attribute sheet: insert a new column before the old "Description" note box ---
$wsAttr = $wb.Worksheets.Item("Attribute")
$wsAttr.Range("L1").EntireColumn.Insert()

# Fill the newly freed / inserted header cells.
# "Places" is typed first (so it lands earlier in the shared-strings table),
# then "MotivationSkills" is typed into the column to its left.
$wsAttr.Range("L1").Value = "Places"
$wsAttr.Range("K1").Value = "MotivationSkills"

# --- LootChance sheet: move the selection, make it no longer the active tab ---
$wsLoot = $wb.Worksheets.Item("LootChance")
$wsLoot.Activate()
$wsLoot.Range("N1").Select()

# --- Attribute sheet becomes the active tab again, with its own selection ---
$wsAttr.Activate()
$wsAttr.Range("M8:O8").Select()
